{"js": "// Load all paragraphs of the document body so we can locate the\n// anchor points by their text content.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// --- Change 1: insert a new \"Tools - ...\" paragraph right after the\n// \"Github - ...\" paragraph (the last bullet of the \"Communication\n// policies, procedures, and tools\" section), i.e. right before the\n// following \"Configuration Management\" heading. ---\nconst toolsText =\n  \"Tools - The project will be using a custom stack consisting of React, \" +\n  \"Django, and Postgres (RPD), where React provides a front end framework \" +\n  \"for building web and mobile applications and allows us to fetch page \" +\n  \"specific javascript allowing for a simple to use SPA. Django provides \" +\n  \"the REST API, user authentication, and serves static files via a \" +\n  \"reverse proxy with NGINX and postgres is a robust relational database\";\n\n// Anchor on the \"Github - ...\" paragraph (the last bullet of the\n// \"Communication policies, procedures, and tools\" section, right before\n// the \"Configuration Management\" heading) and insert the new paragraph\n// right after it, so the new paragraph inherits the plain body style\n// (no pStyle) instead of the following heading's style.\nlet githubParagraph = null;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Github - Website for hosting the git repository\") === 0) {\n    githubParagraph = items[i];\n    break;\n  }\n}\n\nif (githubParagraph) {\n  githubParagraph.insertParagraph(toolsText, Word.InsertLocation.after);\n}\n\n// --- Change 2: change the \"Severity - High\" bullet that sits under the\n// \"Database Structure\" risk-analysis item to \"Severity - Medium\" (the\n// other \"Severity - High\" bullet, under \"User Authentication\", stays\n// untouched). ---\nlet inDatabaseStructureSection = false;\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text.trim();\n  if (text === \"Database Structure\") {\n    inDatabaseStructureSection = true;\n    continue;\n  }\n  if (inDatabaseStructureSection) {\n    if (text === \"Severity - High\") {\n      items[i].insertText(\"Severity - Medium\", Word.InsertLocation.replace);\n      break;\n    }\n    // Any other top-level bullet means we've left the section without\n    // finding the target paragraph; stop looking just in case.\n    if (text === \"User Experience\" || text === \"Queue System Issues\" || text === \"Hosting\") {\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Change 1: insert a new \"Tools - ...\" paragraph right after the\n# \"Github - ...\" paragraph (the last bullet of the \"Communication\n# policies, procedures, and tools\" section), i.e. right before the\n# following \"Configuration Management\" heading. ---\n$toolsText = \"Tools - The project will be using a custom stack consisting of React, Django, and Postgres (RPD), where React provides a front end framework for building web and mobile applications and allows us to fetch page specific javascript allowing for a simple to use SPA. Django provides the REST API, user authentication, and serves static files via a reverse proxy with NGINX and postgres is a robust relational database\"\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -like \"Github - Website for hosting the git repository*\") {\n        $p.Range.InsertParagraphAfter()\n        $newP = $d.Paragraphs.Item($i + 1)\n        $newP.Range.Text = $toolsText\n        break\n    }\n}\n\n# --- Change 2: change the \"Severity - High\" bullet that sits under the\n# \"Database Structure\" risk-analysis item to \"Severity - Medium\" (the\n# other \"Severity - High\" bullet, under \"User Authentication\", stays\n# untouched). ---\n$count = $d.Paragraphs.Count\n$inDatabaseStructureSection = $false\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd(\"`r\", \"`n\", [char]7)\n    if ($t -eq \"Database Structure\") {\n        $inDatabaseStructureSection = $true\n        continue\n    }\n    if ($inDatabaseStructureSection) {\n        if ($t -eq \"Severity - High\") {\n            $p.Range.Text = \"Severity - Medium\"\n            break\n        }\n        if ($t -eq \"User Experience\" -or $t -eq \"Queue System Issues\" -or $t -eq \"Hosting\") {\n            break\n        }\n    }\n}\n"}
